$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the typo in F52: "Imrpove Code" -> "Improve Code" ---
$ws.Range("F52").Value = "Improve Code"

# --- Row 53: 13:30 - 13:45, Source MUX / Improve Code / Add packages ---
$ws.Range("A53").Value = "'6.4.2020"
$ws.Range("B53").Value = 0.5625
$ws.Range("C53").Value = 0.57291666666666663
$ws.Range("D53").Formula = "=C53-B53"
$ws.Range("E53").Value = "Source MUX"
$ws.Range("F53").Value = "Improve Code"
$ws.Range("G53").Value = "Add packages"

# --- Row 54: 13:45 - 14:00, VGA Control / Improve Code / Add packages ---
$ws.Range("A54").Value = "'6.4.2020"
$ws.Range("B54").Value = 0.57291666666666663
$ws.Range("C54").Value = 0.58333333333333337
$ws.Range("D54").Formula = "=C54-B54"
$ws.Range("E54").Value = "VGA Control"
$ws.Range("F54").Value = "Improve Code"
$ws.Range("G54").Value = "Add packages"

# --- Copy formatting from row 51 (a fully/normally formatted data row) down ---
# onto the two new rows, so the new cells pick up the same style indices
# used throughout the table (A: s=2, B/C: s=5, D: s=3, E/F: s=6, G: none),
# after the values have already been entered so the text-vs-date coercion
# for column A has already settled.
$ws.Range("A51:G51").Copy()
$ws.Range("A53:G54").PasteSpecial(-4122)

# --- Move the sheet's selection like a user would after typing the new rows ---
$ws.Range("A55").Select()
